# "Nome do Projeto + business model canvas"
#
# 1. The subtitle paragraph ("Proposta de Projeto em Sistemas de
#    Informação") loses the stray direct formatting on its paragraph
#    mark (Helvetica/colour/size) picked up from the template and is
#    left with just the italic Trebuchet/Times run formatting that its
#    own runs already use.
# 2. A new, centred, bold project-title paragraph ("HistoryCar") is
#    added right after it.

$d = $word.ActiveDocument

# --- 1. Normalise the paragraph-mark formatting of the subtitle paragraph ---

$subtitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "Proposta de Projeto em Sistemas de Informa*") {
        $subtitlePara = $candidate
        break
    }
}

$subtitleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="3C843B55" w14:textId="77777777" w:rsidR="00B87696" w:rsidRPr="008C1789" w:rsidRDefault="00B87696">
<w:pPr>
<w:pStyle w:val="Subttulo"/>
<w:rPr>
<w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Times"/>
<w:i/>
</w:rPr>
</w:pPr>
<w:r w:rsidRPr="008C1789">
<w:rPr>
<w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Times"/>
<w:i/>
</w:rPr>
<w:t>Proposta de Projeto</w:t>
</w:r>
<w:r w:rsidR="00272923">
<w:rPr>
<w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Times"/>
<w:i/>
</w:rPr>
<w:t xml:space="preserve"> em</w:t>
</w:r>
<w:r w:rsidR="00B95D24">
<w:rPr>
<w:rFonts w:ascii="Trebuchet MS" w:hAnsi="Trebuchet MS" w:cs="Times"/>
<w:i/>
</w:rPr>
<w:t xml:space="preserve"> Sistemas de Informação</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$subtitlePara.Range.InsertXML($subtitleXml)

# --- 2. Insert the new project-name paragraph ("HistoryCar") right after it ---

$subtitlePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "Proposta de Projeto em Sistemas de Informa*") {
        $subtitlePara = $candidate
        break
    }
}

$afterRange = $subtitlePara.Range
$afterRange.Collapse(0)
$afterRange.InsertParagraphAfter()

$titlePara = $d.Paragraphs($subtitlePara.Index + 1)

$titleXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:jc w:val="center"/>
<w:rPr>
<w:sz w:val="44"/>
<w:szCs w:val="44"/>
</w:rPr>
</w:pPr>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:rPr>
<w:b/>
<w:bCs/>
<w:sz w:val="44"/>
<w:szCs w:val="44"/>
</w:rPr>
<w:t>HistoryCar</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$titlePara.Range.InsertXML($titleXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
